# This workbook was round-tripped through an unlicensed Aspose.Cells for
# .NET conversion (part of adding .xls upload support / file-type
# filtering to a FileUpload control). Aspose's evaluation mode:
#   1) renames the original worksheet (Hoja1 -> Hoja4), and
#   2) appends a new "Evaluation Warning" worksheet carrying its
#      watermark text, which becomes the active/selected tab.
# Reproduce both effects here.

$wb = $excel.ActiveWorkbook

# 1) Rename the original worksheet.
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Name = "Hoja4"
[void]$ws1.Range("A1").Select()
$ws1.PageSetup.Orientation = 1

# 2) Append the evaluation-warning worksheet right after it.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Evaluation Warning"

$cell = $ws2.Range("A5")
$cell.Value = "Evaluation Only. Created with Aspose.Cells for .NET.Copyright 2003 - 2022 Aspose Pty Ltd."
$cell.Font.Name = "Arial"
$cell.Font.Size = 18
$cell.Font.Bold = $true
$cell.Font.Italic = $true
$cell.Font.Color = 16711680   # RGB(0,0,255) -> blue, BGR-encoded as Excel expects

$ws2.Rows.Item(5).RowHeight = 23.25
$ws2.PageSetup.Orientation = 1

# Aspose leaves the warning sheet as the selected/active tab on open.
$ws2.Activate()
